$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated price/volume figures for the symbol list refresh.
# Each target cell stores a literal text value (e.g. "274.38" or "1.57%"),
# so we force a Text number format before the write and restore the default
# "Normal" style afterwards to avoid leaving a stray number format behind.
$updates = @(
    @{ Cell = "D2"; Value = "274.38" }
    @{ Cell = "E2"; Value = "1.57%" }
    @{ Cell = "D3"; Value = "26.82" }
    @{ Cell = "E3"; Value = "0.44%" }
    @{ Cell = "D4"; Value = "4.921" }
    @{ Cell = "E4"; Value = "4.73%" }
    @{ Cell = "D5"; Value = "0.06343" }
    @{ Cell = "E5"; Value = "3.93%" }
    @{ Cell = "D6"; Value = "6.948" }
    @{ Cell = "E6"; Value = "3.11%" }
    @{ Cell = "D7"; Value = "3.359" }
    @{ Cell = "E7"; Value = "6.04%" }
    @{ Cell = "D8"; Value = "1.431" }
    @{ Cell = "E8"; Value = "60.23%" }
    @{ Cell = "D9"; Value = "0.8884" }
    @{ Cell = "D10"; Value = "0.1471" }
    @{ Cell = "E10"; Value = "3.72%" }
    @{ Cell = "D11"; Value = "0.04966" }
    @{ Cell = "E11"; Value = "-1.36%" }
    @{ Cell = "D12"; Value = "0.07403" }
    @{ Cell = "E12"; Value = "4.32%" }
    @{ Cell = "D13"; Value = "0.03167" }
    @{ Cell = "E13"; Value = "-0.04%" }
    @{ Cell = "D14"; Value = "0.09062" }
    @{ Cell = "E14"; Value = "0.34%" }
    @{ Cell = "D15"; Value = "0.001569" }
    @{ Cell = "E15"; Value = "2.68%" }
    @{ Cell = "D16"; Value = "0.0006339" }
    @{ Cell = "E16"; Value = "4.56%" }
    @{ Cell = "D17"; Value = "0.006034" }
    @{ Cell = "E17"; Value = "-0.89%" }
    @{ Cell = "E18"; Value = "0.64%" }
    @{ Cell = "E19"; Value = "1.66%" }
    @{ Cell = "D20"; Value = "0.3155" }
    @{ Cell = "E20"; Value = "2.22%" }
    @{ Cell = "E21"; Value = "2.82%" }
    @{ Cell = "D22"; Value = "3.908" }
    @{ Cell = "E22"; Value = "1.72%" }
    @{ Cell = "D23"; Value = "0.04343" }
    @{ Cell = "E23"; Value = "2.41%" }
    @{ Cell = "D24"; Value = "0.001177" }
    @{ Cell = "E24"; Value = "-0.76%" }
    @{ Cell = "D25"; Value = "0.003651" }
    @{ Cell = "E25"; Value = "-12.07%" }
    @{ Cell = "D26"; Value = "0.0001203" }
    @{ Cell = "E26"; Value = "0.18%" }
    @{ Cell = "D27"; Value = "0.0001945" }
    @{ Cell = "E27"; Value = "15.77%" }
    @{ Cell = "D40"; Value = "0.04030" }
    @{ Cell = "E40"; Value = "1.91%" }
    @{ Cell = "D41"; Value = "0.006636" }
    @{ Cell = "D42"; Value = "0.1169" }
    @{ Cell = "E42"; Value = "4.66%" }
    @{ Cell = "D43"; Value = "0.002365" }
    @{ Cell = "E43"; Value = "17.68%" }
    @{ Cell = "D44"; Value = "0.01264" }
    @{ Cell = "E44"; Value = "0.61%" }
    @{ Cell = "D45"; Value = "0.00005275" }
    @{ Cell = "E45"; Value = "2.79%" }
    @{ Cell = "E46"; Value = "1,031.62%" }
    @{ Cell = "D47"; Value = "0.02129" }
    @{ Cell = "E47"; Value = "-13.03%" }
    @{ Cell = "E48"; Value = "-0.01%" }
)

foreach ($u in $updates) {
    $c = $ws.Range($u.Cell)
    $c.NumberFormat = "@"
    $c.Value = $u.Value
    $c.Style = "Normal"
}

Write-Output "Updated $($updates.Count) cells"
